$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.118095666666667
$ws.Range("H2").Value = 3.354287
$ws.Range("I2").Value = 0.1244546022965661
$ws.Range("J2").Value = 0.124454602296566
$ws.Range("M2").Value = 0.02507166666666667
$ws.Range("N2").Value = 0.075215
$ws.Range("O2").Value = 0.009392568139045224
$ws.Range("P2").Value = 0.009392568139045224
$ws.Range("Q2").Value = 0.02803252185611112
$ws.Range("R2").Value = 0.252292696705
$ws.Range("S2").Value = 0.001168948332288271
$ws.Range("T2").Value = 0.001168948332288271
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.118095666666667
$ws.Range("H3").Value = 3.354287
$ws.Range("I3").Value = 0.1244546022965661
$ws.Range("J3").Value = 0.124454602296566
$ws.Range("N3").Value = 7.038411000000001
$ws.Range("O3").Value = 0.8789304647757153
$ws.Range("P3").Value = 0.8789304647757155
$ws.Range("Q3").Value = 2.623205613106334
$ws.Range("R3").Value = 23.60885051795701
$ws.Range("S3").Value = 0.1093869414399976
$ws.Range("T3").Value = 0.1093869414399976
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.118095666666667
$ws.Range("H4").Value = 3.354287
$ws.Range("I4").Value = 0.1244546022965661
$ws.Range("J4").Value = 0.124454602296566
$ws.Range("M4").Value = 0.2981003333333334
$ws.Range("N4").Value = 0.894301
$ws.Range("O4").Value = 0.1116769670852394
$ws.Range("P4").Value = 0.1116769670852394
$ws.Range("Q4").Value = 0.333304690931889
$ws.Range("R4").Value = 2.999742218387
$ws.Range("S4").Value = 0.01389871252428017
$ws.Range("T4").Value = 0.01389871252428017
$ws.Range("I5").Value = 0.03090777448944957
$ws.Range("J5").Value = 0.03090777448944957
$ws.Range("M5").Value = 0.02507166666666667
$ws.Range("N5").Value = 0.075215
$ws.Range("O5").Value = 0.009392568139045224
$ws.Range("P5").Value = 0.009392568139045224
$ws.Range("Q5").Value = 0.006961758327222223
$ws.Range("R5").Value = 0.06265582494500001
$ws.Range("S5").Value = 0.0002903033779183988
$ws.Range("T5").Value = 0.0002903033779183988
$ws.Range("I6").Value = 0.03090777448944957
$ws.Range("J6").Value = 0.03090777448944957
$ws.Range("N6").Value = 7.038411000000001
$ws.Range("O6").Value = 0.8789304647757153
$ws.Range("P6").Value = 0.8789304647757155
$ws.Range("S6").Value = 0.02716578459719491
$ws.Range("T6").Value = 0.02716578459719491
$ws.Range("I7").Value = 0.03090777448944957
$ws.Range("J7").Value = 0.03090777448944957
$ws.Range("M7").Value = 0.2981003333333334
$ws.Range("N7").Value = 0.894301
$ws.Range("O7").Value = 0.1116769670852394
$ws.Range("P7").Value = 0.1116769670852394
$ws.Range("Q7").Value = 0.08277481132477779
$ws.Range("R7").Value = 0.744973301923
$ws.Range("S7").Value = 0.003451686514336263
$ws.Range("T7").Value = 0.003451686514336262
$ws.Range("G8").Value = 2.638285666666667
$ws.Range("H8").Value = 7.914857
$ws.Range("I8").Value = 0.293666099582174
$ws.Range("J8").Value = 0.293666099582174
$ws.Range("M8").Value = 0.02507166666666667
$ws.Range("N8").Value = 0.075215
$ws.Range("O8").Value = 0.009392568139045224
$ws.Range("P8").Value = 0.009392568139045224
$ws.Range("Q8").Value = 0.06614621880611113
$ws.Range("R8").Value = 0.5953159692550001
$ws.Range("S8").Value = 0.00275827885045321
$ws.Range("T8").Value = 0.002758278850453209
$ws.Range("G9").Value = 2.638285666666667
$ws.Range("H9").Value = 7.914857
$ws.Range("I9").Value = 0.293666099582174
$ws.Range("J9").Value = 0.293666099582174
$ws.Range("N9").Value = 7.038411000000001
$ws.Range("O9").Value = 0.8789304647757153
$ws.Range("P9").Value = 0.8789304647757155
$ws.Range("Q9").Value = 6.189779619136335
$ws.Range("R9").Value = 55.70801657222701
$ws.Range("S9").Value = 0.2581120813946318
$ws.Range("T9").Value = 0.2581120813946317
$ws.Range("G10").Value = 2.638285666666667
$ws.Range("H10").Value = 7.914857
$ws.Range("I10").Value = 0.293666099582174
$ws.Range("J10").Value = 0.293666099582174
$ws.Range("M10").Value = 0.2981003333333334
$ws.Range("N10").Value = 0.894301
$ws.Range("O10").Value = 0.1116769670852394
$ws.Range("P10").Value = 0.1116769670852394
$ws.Range("Q10").Value = 0.7864738366618891
$ws.Range("R10").Value = 7.078264529957001
$ws.Range("S10").Value = 0.03279573933708909
$ws.Range("T10").Value = 0.03279573933708908
$ws.Range("G11").Value = 1.628177666666667
$ws.Range("H11").Value = 4.884533
$ws.Range("I11").Value = 0.1812315439673029
$ws.Range("J11").Value = 0.1812315439673029
$ws.Range("M11").Value = 0.02507166666666667
$ws.Range("N11").Value = 0.075215
$ws.Range("O11").Value = 0.009392568139045224
$ws.Range("P11").Value = 0.009392568139045224
$ws.Range("Q11").Value = 0.04082112773277778
$ws.Range("R11").Value = 0.367390149595
$ws.Range("S11").Value = 0.001702229625657263
$ws.Range("T11").Value = 0.001702229625657263
$ws.Range("G12").Value = 1.628177666666667
$ws.Range("H12").Value = 4.884533
$ws.Range("I12").Value = 0.1812315439673029
$ws.Range("J12").Value = 0.1812315439673029
$ws.Range("N12").Value = 7.038411000000001
$ws.Range("O12").Value = 0.8789304647757153
$ws.Range("P12").Value = 0.8789304647757155
$ws.Range("Q12").Value = 3.819927866340334
$ws.Range("R12").Value = 34.379350797063
$ws.Range("S12").Value = 0.1592899251712021
$ws.Range("T12").Value = 0.1592899251712021
$ws.Range("G13").Value = 1.628177666666667
$ws.Range("H13").Value = 4.884533
$ws.Range("I13").Value = 0.1812315439673029
$ws.Range("J13").Value = 0.1812315439673029
$ws.Range("M13").Value = 0.2981003333333334
$ws.Range("N13").Value = 0.894301
$ws.Range("O13").Value = 0.1116769670852394
$ws.Range("P13").Value = 0.1116769670852394
$ws.Range("Q13").Value = 0.4853603051592223
$ws.Range("R13").Value = 4.368242746433
$ws.Range("S13").Value = 0.02023938917044361
$ws.Range("T13").Value = 0.02023938917044361
$ws.Range("G14").Value = 1.720859666666667
$ws.Range("H14").Value = 5.162579
$ws.Range("I14").Value = 0.1915479254666055
$ws.Range("J14").Value = 0.1915479254666054
$ws.Range("M14").Value = 0.02507166666666667
$ws.Range("N14").Value = 0.075215
$ws.Range("O14").Value = 0.009392568139045224
$ws.Range("P14").Value = 0.009392568139045224
$ws.Range("Q14").Value = 0.04314481994277778
$ws.Range("R14").Value = 0.388303379485
$ws.Range("S14").Value = 0.001799126941837848
$ws.Range("T14").Value = 0.001799126941837847
$ws.Range("G15").Value = 1.720859666666667
$ws.Range("H15").Value = 5.162579
$ws.Range("I15").Value = 0.1915479254666055
$ws.Range("J15").Value = 0.1915479254666054
$ws.Range("N15").Value = 7.038411000000001
$ws.Range("O15").Value = 0.8789304647757153
$ws.Range("P15").Value = 0.8789304647757155
$ws.Range("Q15").Value = 4.037372535774334
$ws.Range("R15").Value = 36.33635282196901
$ws.Range("S15").Value = 0.1683573071571876
$ws.Range("T15").Value = 0.1683573071571876
$ws.Range("G16").Value = 1.720859666666667
$ws.Range("H16").Value = 5.162579
$ws.Range("I16").Value = 0.1915479254666055
$ws.Range("J16").Value = 0.1915479254666054
$ws.Range("M16").Value = 0.2981003333333334
$ws.Range("N16").Value = 0.894301
$ws.Range("O16").Value = 0.1116769670852394
$ws.Range("P16").Value = 0.1116769670852394
$ws.Range("Q16").Value = 0.5129888402532222
$ws.Range("R16").Value = 4.616899562279
$ws.Range("S16").Value = 0.02139149136757999
$ws.Range("T16").Value = 0.02139149136757999
$ws.Range("G17").Value = 1.600871
$ws.Range("H17").Value = 4.802613
$ws.Range("I17").Value = 0.1781920541979019
$ws.Range("J17").Value = 0.1781920541979019
$ws.Range("M17").Value = 0.02507166666666667
$ws.Range("N17").Value = 0.075215
$ws.Range("O17").Value = 0.009392568139045224
$ws.Range("P17").Value = 0.009392568139045224
$ws.Range("Q17").Value = 0.04013650408833334
$ws.Range("R17").Value = 0.361228536795
$ws.Range("S17").Value = 0.001673681010890233
$ws.Range("T17").Value = 0.001673681010890233
$ws.Range("G18").Value = 1.600871
$ws.Range("H18").Value = 4.802613
$ws.Range("I18").Value = 0.1781920541979019
$ws.Range("J18").Value = 0.1781920541979019
$ws.Range("N18").Value = 7.038411000000001
$ws.Range("O18").Value = 0.8789304647757153
$ws.Range("P18").Value = 0.8789304647757155
$ws.Range("Q18").Value = 3.755862685327
$ws.Range("R18").Value = 33.802764167943
$ws.Range("S18").Value = 0.1566184250155014
$ws.Range("T18").Value = 0.1566184250155014
$ws.Range("G19").Value = 1.600871
$ws.Range("H19").Value = 4.802613
$ws.Range("I19").Value = 0.1781920541979019
$ws.Range("J19").Value = 0.1781920541979019
$ws.Range("M19").Value = 0.2981003333333334
$ws.Range("N19").Value = 0.894301
$ws.Range("O19").Value = 0.1116769670852394
$ws.Range("P19").Value = 0.1116769670852394
$ws.Range("Q19").Value = 0.4772201787236667
$ws.Range("R19").Value = 4.294981608513
$ws.Range("S19").Value = 0.01989994817151029
$ws.Range("T19").Value = 0.01989994817151029
